# Apply the Fri Jun 14 2024 cryptos-list refresh (GitHub Actions data pull).
# Updates price (D) / 1h-volume-change (E) columns, and for the handful of rows
# whose ranking flipped, the coin name (B) + link (C) as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '66.973.17'
$ws.Range('E2').Value = '  -0.89%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.490.89'
$ws.Range('E3').Value = '  -0.90%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.00%  '

# Row 5: BNB
$ws.Range('D5').Value = '''606.43'
$ws.Range('E5').Value = '  -0.25%  '

# Row 6: Solana
$ws.Range('D6').Value = '''145.38'
$ws.Range('E6').Value = '  -2.84%  '

# Row 7: LidoStakedEther
$ws.Range('D7').Value = '3.489.86'
$ws.Range('E7').Value = '  -0.78%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.01%  '

# Row 9: XRP
$ws.Range('D9').Value = '''0.478'
$ws.Range('E9').Value = '  -1.76%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -2.07%  '

# Row 11: Toncoin
$ws.Range('E11').Value = '  +4.53%  '

# Row 12: Cardano
$ws.Range('D12').Value = '''0.419'
$ws.Range('E12').Value = '  -2.49%  '

# Row 13: ShibaInu
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '''0.0000213'
$ws.Range('E13').Value = '  -1.32%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.076.85'
$ws.Range('E14').Value = '  -0.87%  '

# Row 15: Avalanche
$ws.Range('D15').Value = '''31.11'
$ws.Range('E15').Value = '  -3.08%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '3.494.70'
$ws.Range('E16').Value = '  -0.47%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '66.986.59'
$ws.Range('E17').Value = '  -1.00%  '

# Row 18: TRON
$ws.Range('E18').Value = '  -0.29%  '

# Row 19: Uniswap
$ws.Range('D19').Value = '''10.71'
$ws.Range('E19').Value = '  +7.13%  '

# Row 20: Polkadot
$ws.Range('D20').Value = '''6.30'
$ws.Range('E20').Value = '  -3.28%  '

# Row 21: Chainlink
$ws.Range('D21').Value = '''15.32'
$ws.Range('E21').Value = '  -1.49%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''431.24'
$ws.Range('E22').Value = '  -4.46%  '

# Row 23: Polygon
$ws.Range('D23').Value = '''0.603'
$ws.Range('E23').Value = '  -3.80%  '

# Row 24: Litecoin
$ws.Range('D24').Value = '''79.52'
$ws.Range('E24').Value = '  +0.68%  '

# Row 25: Dai
$ws.Range('E25').Value = '  +0.10%  '

# Row 26: WrappedeETH
$ws.Range('D26').Value = '3.624.23'

# Row 27: PEPE
$ws.Range('E27').Value = '  -5.43%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range('D28').Value = '''9.74'
$ws.Range('E28').Value = '  -2.14%  '

# Row 29: RenderToken
$ws.Range('D29').Value = '''8.13'
$ws.Range('E29').Value = '  -5.40%  '

# Row 30: PancakeSwap
$ws.Range('E30').Value = '  -0.66%  '

# Row 31: Fetch.AI
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '''1.55'
$ws.Range('E31').Value = '  -6.12%  '

# Row 32: Binance-PegBSC-USD
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.02%  '

# Row 33: Kaspa
$ws.Range('D33').Value = '''0.166'
$ws.Range('E33').Value = '  -1.90%  '

# Row 34: EthereumClassic
$ws.Range('D34').Value = '''25.37'
$ws.Range('E34').Value = '  -1.33%  '

# Row 35: ImmutableX
$ws.Range('D35').Value = '''1.78'
$ws.Range('E35').Value = '  -3.43%  '

# Row 37: Aptos
$ws.Range('D37').Value = '''7.93'
$ws.Range('E37').Value = '  -1.07%  '

# Row 38: NEARProtocol
$ws.Range('D38').Value = '''5.71'
$ws.Range('E38').Value = '  -8.32%  '

# Row 39: FirstDigitalUSD
$ws.Range('D39').Value = '''0.998'
$ws.Range('E39').Value = '  +0.01%  '

# Row 40: Monero
$ws.Range('D40').Value = '''173.89'
$ws.Range('E40').Value = '  -1.17%  '

# Row 41: Hedera
$ws.Range('D41').Value = '''0.0894'
$ws.Range('E41').Value = '  -0.68%  '

# Row 42: Filecoin
$ws.Range('D42').Value = '''5.34'
$ws.Range('E42').Value = '  -1.60%  '

# Row 43: Stacks
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.99'
$ws.Range('E43').Value = '  -12.62%  '

# Row 44: Mantle
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '''0.892'
$ws.Range('E44').Value = '  -1.01%  '

# Row 45: OKB
$ws.Range('D45').Value = '''46.35'
$ws.Range('E45').Value = '  -1.18%  '

# Row 46: InjectiveProtocol
$ws.Range('D46').Value = '''27.65'
$ws.Range('E46').Value = '  -10.56%  '

# Row 47: ONDO
$ws.Range('E47').Value = '  -5.89%  '

# Row 48: Cosmos
$ws.Range('D48').Value = '''7.30'
$ws.Range('E48').Value = '  -4.07%  '

# Row 49: dogwifhat
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '''2.39'
$ws.Range('E49').Value = '  -3.86%  '

# Row 50: SuiNetwork
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').Value = '''0.982'
$ws.Range('E50').Value = '  -1.74%  '

# Row 51: TheGraph
$ws.Range('D51').Value = '''0.244'
$ws.Range('E51').Value = '  -2.79%  '
